# Review_442.docx update (per commit: daily-paper switched from the
# "Training Large Language Models to Reason in a Continuous Latent Space"
# review to the "Critical Tokens Matter" review).
#
# The document has 7 paragraphs. Paragraph 1 holds two runs of text (the date
# line and, after a manual line break, the paper title); paragraphs 2-6 are the
# Hebrew review body; paragraph 7 is the arxiv link. All eight text changes from
# the diff are applied below. Most are applied with Find/Replace; the body
# paragraph whose replacement text contains literal straight double-quotes is
# instead written directly via Range.Text so Word's smart-quote autocorrect
# does not turn them into curly quotes.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# 1) Date line: 20.04.25 -> 18.04.25
Replace-Text " המאמר היומי של מייק: 20.04.25" " המאמר היומי של מייק: 18.04.25"

# 2) Paper title (second half of paragraph 1, after the <w:br/>)
Replace-Text "Training Large Language Models to Reason in a Continuous Latent Space" "Critical Tokens Matter: Token-Level Contrastive Estimation Enhances LLM’s Reasoning Capability"

# 3) Paragraph 2 - new opening summary
Replace-Text "המאמר מציג רעיון חדשני ומתבקש (לעניות דעתי) לשיפור תהליכי הנמקה (reasoning) של מודל שפה. כמו שאתם בטח יודעים אנו גורמים למודלי שפה לחשוב על ידי הכנסה לפרומפט ביטוים כמו ״think step by step״ או טוקנים מיוחדים של חשיבה כמו <think> וכדומה. זה גורם למודל ״לפלוט״ את שרשרת הנמקה בצורה של טוקנים, כלומר של טקסט. היתרון בגישות אלו שאנו יכולים לנתח את שרשרת החשיבה של מודל ולשפר אותה כי אנו רואים אותה כטקסט." "מאמר די מעניין העוסק בשיפור יכולות הנמקה של מודלי שפה בשאלות שיש להם תשובות חד משמעיות (כגון שאלות מתמטיות ושאלות קוד שניתן לבחון את נכונות הפתרון באמצעות סט מקיף של טסטים). המאמר מגדיר מושג טוקן קריטי (critical token) שהוא למעשה מהווה סוג של סימן האם המודל הולך לתת תשובה נכונה או לא נכונה לשאלה. "

# 4) Paragraph 3 - contains literal straight double-quotes in the replacement,
#    so it is written via Range.Text (Find/Replace would smart-quote it).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.Contains("DeepSeek")) {
        $targetPara = $candidate
        break
    }
}
if ($targetPara -eq $null) {
    throw "could not locate the DeepSeek paragraph to replace"
}
$targetPara.Range.Text = "המחברים שמו לב כי שבתוך מסלולי הנמקה(reasoning) שגויים, ישנם טוקנים שהם כמעט בוודאות מובילים לתוצאות שגויות. טוקנים אלה משבשים את הרצף הלוגי, מעוותים קשרים או מכניסים שגיאות חישוביות, וכך משפיעים באופן משמעותי על התוצאה הסופית. בשונה מטוקנים אחרים שעשויים להשפיע בצורה לא משמעותית על תהליך האינפרנס, `"הטוקנים הקריטיים`" האלו מהוות סוג של נקודות כשל. זיהוי הטוקנים הללו הוא חיוני, משום שלעיתים קרובות הימנעות מהם או תיקונם יכולה להוביל לתוצאה נכונה – גם בתוך מסלול הסקה שגוי."

# 5) Paragraph 4 - critical-token identification method
Replace-Text "המאמר שנסקור היום עושה צעד נוסף בכיוון הזה. הרי מודלי שפה לא חייבים לחשוב בשפות שאנו, בני אדם, מבינים, נכון? בשביל כך יש להם את מחרב הייצוג שלהם, כלומר המרחב הלטנטי. הרי מודל שפה לא חושב באמצעות מילים ובמשפטים כמונו אלא פועל במרחב וקטורי שכל וקטור ייצוג של טוקן. אז המחברים אמרו את הדבר הבא: בוא נחליף שרשרת הנמקה בשרשראות הנמקה לטנטיות (וקטוריות) ללא תרגום לשפה האנושית. אז המודל מאומן להחליף שרשראות הנמקה בשפה טבעית בסדרה של וקטורים. " "המאמר מציעה שיטה לזיהוי של טוקנים קריטיים. טוקן מזוהה כקריטי עם כל מסלולי ההנמקה המתחילים ממנו מסתיימים בתשובה שגויה ועבור כל הטוקנים שבאים 95% מהמסלולים המתחילים מהם מסתיימים בתשובה לא נכונה. שימו לב שיש טוקנים המופיעים בטקסט במיקומים שהם אחרי הטוקן הקריטי שלא כל מסלולי הנמקה שלהם מכילים את הטוקן הקריטי, כך לא מן הנמנע שיש בינם מסלולים המסתיימים בתשובה נכונה. המחברים ביצועו כמה בדיקות כדי לוודא שהטוקנים שזוהו בצורה באמת טוקנים קריטים."

# 6) Paragraph 5 - RLHF alignment discussion
Replace-Text "וזה בדיוק מה שנעשה באימון המודל. המחברים מאמנים מודל לפלוט וקטורים במקום עבור כמה שלבי הנמקה ראשונים. כלומר המודל מאומן (בו זמנית) להחליף שלבים 1-3 או 1-6 של שרשרת הנמקה בוקטורים. כלומר המודל מתחיל מהמשטר הלטנטי (latent mode) שהמחשבות שלו הם הוקטורים וממשיך במשטר שפתי (language mode) שבו הפלט הוא שפה טבעית. כמובן שיש טוקן שמפריד בין משטרים אלו כלומר <eot>." "לאחר מכן המאמר מפתח שיטת RLHF ליישור מודל שפה שבמרכזה מזעור של הנראות של הטוקנים הקריטיים (כי הם מובילים לשגיאות). בשביל כך המאמר מציע לאמן שני מודלים (עם פיינטיון) - אחד שמגנרט תשבות נכונות והשני מגנרט תשובות לא נכונות (שמעתם נכון). "

# 7) Paragraph 6 - DPO formula discussion (replaces the old closing remark)
Replace-Text "מאמר עם כיוון מאוד מעניין שאני צופה לו עתיד גדול." "לאחר מכם המחברים מנסחים דרך לשערוך הנראות של האם הטוקן הוא קריטי בהינתן הפרומפט וטוקני התשובה לפניו. הנוסחה היא הפרש ממשוקל של נראויות הטוקנים (מותנים) מהמודל של התשובות הנכונות לבין זה של המודל של התשובות השגויות. שערוך זה מקבל ערך נמוך עבור התשובה הנכונה וערך גבוה עבור התשובה הלא נכונה. בשלב האחרון המודל עובר פיינטיון עם DPO שזה קיצור של Direct Preference Optimization כאשר על הזוגות של שאלות עם התשובות הנכונות והשגויות. כדי למזער את הסיכוי להופעת טוקן קריטי המאמר משנה את האיבר המכיל נראות של תשובה שגויה בנוסחה העיקרית של DPO על ידי הכפלה על הנראות השלילית של טוקן להיות קריטי. שימו לב שמכיוון שההתחשבות בנראות מתרחשת ברמה של טוקן ה-DPO במאמר עובר להיות token-level ולא sample-level כמו במאמר המקורי."

# 8) arxiv link
Replace-Text "https://arxiv.org/pdf/2412.06769" "https://arxiv.org/abs/2411.19943"

